$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("建物")

# The "property_category" column (I) for every data row in the 建物
# (building) sheet was tagged "land" by mistake; correct it to "building".
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    if ($cell.Value2 -eq "land") {
        $cell.Value = "building"
    }
}
